$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.63119999999999
$ws.Range("A7").Value = -20.24929999999998
$ws.Range("A16").Value = -22.01110000000001
$ws.Range("A28").Value = -21.92710000000001
$ws.Range("A29").Value = -21.18869999999998
$ws.Range("A32").Value = -21.1636
$ws.Range("A40").Value = -20.2806
$ws.Range("A52").Value = -22.1807
$ws.Range("A57").Value = -22.58670000000003
$ws.Range("A66").Value = -21.44349999999999
$ws.Range("A100").Value = -22.02480000000001
